# Update "想去人数" (F) and "最低票价" (G) figures across the workbook to the
# newly scraped values (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-CellValues {
    param([string]$SheetName, [hashtable]$Updates)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($ref in $Updates.Keys) {
        $ws.Range($ref).Value = $Updates[$ref]
    }
}

# Sheet "展览" (Exhibitions)
Set-CellValues "展览" @{
    "F2"  = 4606
    "F3"  = 2701
    "F5"  = 2701
    "F10" = 727
    "F12" = 179
    "F13" = 380
    "F14" = 558
    "F15" = 288
    "F18" = 508
    "G18" = 29.9
    "F21" = 634
    "F23" = 134
    "F25" = 486
    "F27" = 1394
    "F28" = 287
    "F29" = 36
    "F30" = 1374
    "F31" = 2236
    "F32" = 359
    "F34" = 587
    "F36" = 47
    "F38" = 752
    "F39" = 1429
    "F40" = 179
    "F42" = 474
    "F43" = 70
}

# Sheet "演出" (Performances)
Set-CellValues "演出" @{
    "F4"  = 82
    "F13" = 14
}

# Sheet "全部类型" (All types)
Set-CellValues "全部类型" @{
    "F2"  = 4606
    "F3"  = 2701
    "F4"  = 2701
    "F8"  = 727
    "F10" = 179
    "F11" = 380
    "F12" = 560
    "F13" = 288
    "F16" = 508
    "G16" = 29.9
    "F18" = 634
    "F20" = 134
    "F21" = 82
    "F25" = 486
    "F27" = 1394
    "F28" = 287
    "F29" = 36
    "F32" = 2236
    "F33" = 359
    "F38" = 587
    "F40" = 47
    "F42" = 752
    "F43" = 1429
    "F45" = 179
    "F46" = 474
    "F47" = 70
    "F49" = 14
}
